$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.036.68'
$ws.Range("E2").Value = '  +1.62%  '

$ws.Range("D3").Value = '3.352.39'
$ws.Range("E3").Value = '  +2.97%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.93'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.47%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '173.11'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.85%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.592'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.39%  '

$ws.Range("D8").Value = '3.351.72'
$ws.Range("E8").Value = '  +3.04%  '

$ws.Range("E9").Value = '  -0.13%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.606'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.72%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.20'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.73%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.133'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.87%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000256'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.01'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.70%  '

$ws.Range("D15").Value = '3.896.42'
$ws.Range("E15").Value = '  +3.55%  '

$ws.Range("D16").Value = '3.357.18'
$ws.Range("E16").Value = '  +3.42%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.117'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").Value = '63.857.67'
$ws.Range("E19").Value = '  +1.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.21'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.960'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.96%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.91'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.25'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.47'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.86%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '81.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.84%  '

$ws.Range("E26").Value = '  -0.41%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.17'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.70%  '

$ws.Range("E28").Value = '  +2.52%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '11.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.70%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.48%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '28.79'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '626.31'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.41'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.99%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '11.17'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.11%  '

$ws.Range("E35").Value = '  +0.20%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '57.86'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '

$ws.Range("E37").Value = '  +0.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '36.03'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.378'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.01%  '

$ws.Range("D40").Value = '0.0₃0730'
$ws.Range("E40").Value = '  +10.28%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.26%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.65'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +8.56%  '

$ws.Range("D43").Value = '2.941.29'
$ws.Range("E43").Value = '  +0.45%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.123'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.75%  '

$ws.Range("E45").Value = '  +6.83%  '

$ws.Range("E46").Value = '  +2.09%  '

$ws.Range("E47").Value = '  +0.49%  '

$ws.Range("E48").Value = '  -3.35%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.00'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.66%  '

$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '135.64'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.15%  '
